# Update ins.xlsx for January with new data:
# Clear the now-stale instruction rows (4-9) while keeping their cell
# styles intact, and move the active selection to B8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the data (values/types) from rows 4 through 9, columns A:F,
# leaving the formatting (styles) untouched.
$ws.Range("A4:F9").ClearContents()

# Reflect the new selection/active cell saved with the sheet.
$ws.Range("B8").Select()
